$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 27: "Fire Panel" / "Type" / "Fire" collapses into "Fire" row, and the
# "Call everyone..." action text (previously split onto its own row 28 in
# columns I/J) now lives directly on row 27 in columns G/H.
# ---------------------------------------------------------------------------
$ws.Range("F27").Value = "Fire"
$ws.Range("G27").Value = "action"

$callEveryoneText = "Call everyone off the following: <li>Fire Department</li> <li> DPS</li> <li>O&M to have an electrician turn off the alarm </li> <li>Area Director & Building Manager</li> <li>Create a work order for HSM</li> <li> Call Landry </li> <li> Area if necessary to warn the residents</li>"
$ws.Range("H27").Value = $callEveryoneText
$ws.Range("H27").Font.Name = "Arial"
$ws.Range("H27").Font.Size = 8
$ws.Range("H27").Font.Color = 0

$ws.Range("I27").ClearContents()

# ---------------------------------------------------------------------------
# Old row 28 (I28/J28 leftover cells) is no longer needed in that shape --
# clear it out completely (contents + formatting) before repurposing the row
# for new Plumbing data below.
# ---------------------------------------------------------------------------
$ws.Rows.Item(28).Clear()

# ---------------------------------------------------------------------------
# New row 28: Plumbing / Water Shutdown
# ---------------------------------------------------------------------------
$ws.Range("A28").Value = "what time is it"
$ws.Range("B28").Value = "After hours"
$ws.Range("C28").Value = "job type"
$ws.Range("D28").Value = "Plumbing"
$ws.Range("E28").Value = "problem"
$ws.Range("F28").Value = "Water Shutdown"
$ws.Range("G28").Value = "action"
$ws.Range("H28").Value = "Priority B <br> 1 – 8 Hours <br> <li>If <strong>scheduled:</strong> Inform CSC </li><li> If <strong>Emergency shutdown: </strong> Inform resident off the emergency due to which shutdown was done."

# ---------------------------------------------------------------------------
# New row 29: Plumbing / Toilet / Overflowing
# ---------------------------------------------------------------------------
$toiletActionText = "<li> <strong>Priority A: </strong> Page GMT </li><li><strong>Priority B: </strong><li>If water is dripping wait until morning</li><li> If there is stream of water, page gmts</li> </li> "

$ws.Range("A29").Value = "what time is it"
$ws.Range("B29").Value = "After hours"
$ws.Range("C29").Value = "job type"
$ws.Range("D29").Value = "Plumbing"
$ws.Range("E29").Value = "problem"
$ws.Range("F29").Value = "Toilet"
$ws.Range("G29").Value = "Type"
$ws.Range("H29").Value = "Overflowing"
$ws.Range("I29").Value = "action"
$ws.Range("J29").Value = $toiletActionText
$ws.Range("J29").Font.Name = "Arial"
$ws.Range("J29").Font.Size = 8
$ws.Range("J29").Font.Color = 0

# ---------------------------------------------------------------------------
# New row 30: Plumbing / Toilet / Leaking
# ---------------------------------------------------------------------------
$ws.Range("A30").Value = "what time is it"
$ws.Range("B30").Value = "After hours"
$ws.Range("C30").Value = "job type"
$ws.Range("D30").Value = "Plumbing"
$ws.Range("E30").Value = "problem"
$ws.Range("F30").Value = "Toilet"
$ws.Range("G30").Value = "Type"
$ws.Range("H30").Value = "Leaking"
$ws.Range("I30").Value = "action"
$ws.Range("J30").Value = $toiletActionText
$ws.Range("J30").Font.Name = "Arial"
$ws.Range("J30").Font.Size = 8
$ws.Range("J30").Font.Color = 0

# ---------------------------------------------------------------------------
# Update the view: active cell / selection moves to J29, and scroll so that
# D19 is the top-left visible cell (best-effort, mirrors the saved view).
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 19
$win.ScrollColumn = 4
$ws.Range("J29").Select()
